$data = @(
    ,@("header", 14, "Week 3")
    ,@("data", 15, "Tolosa", 1, 2, "Hallacy")
    ,@("data", 16, "Nagle", 3, 4, "Stichler")
    ,@("data", 17, "Rich", 5, 6, "Yamaoka")
    ,@("data", 18, "Pitton", 7, 8, "Netter")
    ,@("data", 19, "Walker", 9, 10, "Nishida")
    ,@("header", 20, "Week 4")
    ,@("data", 21, "Hallacy", 1, 2, "Rich")
    ,@("data", 22, "Pitton", 3, 4, "Tolosa")
    ,@("data", 23, "Nishida", 5, 6, "Yamaoka")
    ,@("data", 24, "Nagle", 7, 8, "Netter")
    ,@("data", 25, "Stichler", 9, 10, "Walker")
    ,@("header", 26, "Week 5")
    ,@("data", 27, "Nishida", 1, 2, "Hallacy")
    ,@("data", 28, "Walker", 3, 4, "Nagle")
    ,@("data", 29, "Yamaoka", 5, 6, "Stichler")
    ,@("data", 30, "Netter", 7, 8, "Tolosa")
    ,@("data", 31, "Rich", 9, 10, "Pitton")
    ,@("header", 32, "Week 6")
    ,@("data", 33, "Stichler", 1, 2, "Hallacy")
    ,@("data", 34, "Tolosa", 3, 4, "Rich")
    ,@("data", 35, "Nishida", 5, 6, "Pitton")
    ,@("data", 36, "Walker", 7, 8, "Netter")
    ,@("data", 37, "Nagle", 9, 10, "Yamaoka")
    ,@("header", 38, "Week 7")
    ,@("data", 39, "Hallacy", 1, 2, "Nagle")
    ,@("data", 40, "Yamaoka", 3, 4, "Walker")
    ,@("data", 41, "Stichler", 5, 6, "Pitton")
    ,@("data", 42, "Netter", 7, 8, "Rich")
    ,@("data", 43, "Nishida", 9, 10, "Tolosa")
    ,@("header", 44, "Week 8")
    ,@("data", 45, "Walker", 1, 2, "Hallacy")
    ,@("data", 46, "Rich", 3, 4, "Nishida")
    ,@("data", 47, "Nagle", 5, 6, "Pitton")
    ,@("data", 48, "Yamaoka", 7, 8, "Netter")
    ,@("data", 49, "Tolosa", 9, 10, "Stichler")
    ,@("header", 50, "Week 9")
    ,@("data", 51, "Hallacy", 1, 2, "Yamaoka")
    ,@("data", 52, "Stichler", 3, 4, "Rich")
    ,@("data", 53, "Nagle", 5, 6, "Tolosa")
    ,@("data", 54, "Netter", 7, 8, "Nishida")
    ,@("data", 55, "Pitton", 9, 10, "Walker")
    ,@("header", 56, "Week 10")
    ,@("data", 57, "Hallacy", 1, 2, "Netter")
    ,@("data", 58, "Yamaoka", 3, 4, "Pitton")
    ,@("data", 59, "Walker", 5, 6, "Tolosa")
    ,@("data", 60, "Nishida", 7, 8, "Stichler")
    ,@("data", 61, "Rich", 9, 10, "Nagle")
    ,@("header", 62, "Week 11")
    ,@("data", 63, "Hallacy", 1, 2, "Pitton")
    ,@("data", 64, "Nagle", 3, 4, "Nishida")
    ,@("data", 65, "Walker", 5, 6, "Rich")
    ,@("data", 66, "Netter", 7, 8, "Stichler")
    ,@("data", 67, "Tolosa", 9, 10, "Yamaoka")
    ,@("header", 68, "Week 12")
    ,@("data", 69, "Hallacy", 1, 2, "Tolosa")
    ,@("data", 70, "Stichler", 3, 4, "Nagle")
    ,@("data", 71, "Yamaoka", 5, 6, "Rich")
    ,@("data", 72, "Pitton", 7, 8, "Netter")
    ,@("data", 73, "Nishida", 9, 10, "Walker")
    ,@("header", 74, "Week 13")
    ,@("data", 75, "Rich", 1, 2, "Hallacy")
    ,@("data", 76, "Tolosa", 3, 4, "Pitton")
    ,@("data", 77, "Yamaoka", 5, 6, "Nishida")
    ,@("data", 78, "Netter", 7, 8, "Nagle")
    ,@("data", 79, "Walker", 9, 10, "Stichler")
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewData")

foreach ($row in $data) {
    $kind = $row[0]
    $r = $row[1]
    if ($kind -eq "header") {
        $text = $row[2]
        $rng = $ws.Range("A$r`:D$r")
        $ws.Range("A$r").Value = $text
        $rng.HorizontalAlignment = -4108
        $rng.Merge()
    } else {
        $a = $row[2]
        $b = $row[3]
        $c = $row[4]
        $d = $row[5]
        $ws.Range("A$r").Value = $a
        $ws.Range("B$r").Value = $b
        $ws.Range("C$r").Value = $c
        $ws.Range("D$r").Value = $d
    }
}

# Update sheet views: NewData becomes the active/selected tab,
# scrolled to show rows around 65-79, with B75:C79 selected.
$ws.Activate()
$ws.Range("B75:C79").Select()

